$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.732.24"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "2.223.25"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.60"
$ws.Range("E5").Value = "  +7.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").Value = "  -1.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.83"
$ws.Range("E7").Value = "  +3.26%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"
$ws.Range("E9").Value = "  +6.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.03"
$ws.Range("E10").Value = "  +11.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0969"
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.31"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.19"
$ws.Range("E13").Value = "  +6.84%  "
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "2.553.87"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.97"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.865"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "2.226.00"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").Value = "41.734.29"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.23"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.89"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.23"
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("E24").Value = "  +6.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.04"
$ws.Range("E25").Value = "  +10.01%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.53"
$ws.Range("E27").Value = "  +5.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.71"
$ws.Range("E28").Value = "  +6.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.17"
$ws.Range("E29").Value = "  +0.75%  "
$ws.Range("E30").Value = "  -4.43%  "
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.125"
$ws.Range("E32").Value = "  +3.19%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.63"
$ws.Range("E33").Value = "  +5.27%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.125"
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.72"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.61"
$ws.Range("E37").Value = "  +16.40%  "
$ws.Range("E38").Value = "  +9.51%  "
$ws.Range("E39").Value = "  +7.29%  "
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.94"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "66.76"
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.11"
$ws.Range("E43").Value = "  +19.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.202"
$ws.Range("E44").Value = "  +5.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.87"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.78"
$ws.Range("E46").Value = "  -4.48%  "
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("E48").Value = "  +2.67%  "
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("E50").Value = "  +5.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.19"
$ws.Range("E51").Value = "  +1.04%  "
